# Update gh-pages to output generated at 456a3b4
# Increments the "想去人数" (want-to-go count, column F) for several events
# across the 展览 / 演出 / 本地生活 sheets, and mirrors the same updates in
# the aggregated 全部类型 sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 26
$ws1.Range("F10").Value = 49
$ws1.Range("F15").Value = 371
$ws1.Range("F23").Value = 265
$ws1.Range("F31").Value = 152

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 530
$ws2.Range("F6").Value = 530
$ws2.Range("F10").Value = 12
$ws2.Range("F16").Value = 713
$ws2.Range("F18").Value = 22

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1741
$ws3.Range("F6").Value = 2209

# Sheet "全部类型" (All Types) - aggregated view of all the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1741
$ws4.Range("F4").Value = 2209
$ws4.Range("F18").Value = 26
$ws4.Range("F21").Value = 530
$ws4.Range("F23").Value = 49
$ws4.Range("F28").Value = 371
$ws4.Range("F35").Value = 265
$ws4.Range("F43").Value = 22
$ws4.Range("F48").Value = 152
